# "Loan RBI, Variable Instalments"
#
# The "Repayment schedule" sheet gets a brand-new (blank) column inserted
# just before the existing "Late" column (column N), pushing the old
# N/O/P columns ("Late", "heading"/Disbursement, "Outstanding") one slot
# to the right (O/P/Q). The sheet also becomes the active tab/selection
# of the workbook (it was "Transactions" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet/tab (was "Transactions").
$ws.Activate() | Out-Null

# Insert a new blank column before column N (14) -- shifts
# Late / heading / Outstanding from N,O,P to O,P,Q.
$ws.Columns.Item(14).Insert() | Out-Null

# New column inherits the width of its left neighbour ("In Advance", col M).
$ws.Columns.Item(14).ColumnWidth = 9.8

# Update the on-screen selection for the now-active sheet.
$ws.Range("O6").Select() | Out-Null
